$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting C:P to D:Q.
$ws.Range("C1").EntireColumn.Insert()

# The new column should inherit column B's width, same as Excel normally does.
$ws.Range("C1").ColumnWidth = $ws.Range("B1").ColumnWidth

# Rename old "Mã đại lý" header (now in B7) to the "auto-generated code" label,
# and give the newly inserted column C its own "manually entered code" header.
$ws.Range("B7").Value = "Mã đại lý (tự sinh)"
$ws.Range("C7").Value = "Mã đại lý (tự nhập)"

# Placeholder row: fill in the new column's merge-field token.
$ws.Range("C9").Value = "{{ReportStoreGenerals.Stores.CodeDraft}}"

$ws.Range("E12").Select()
